# "Actualizar" refresh macro: re-stamp the "Ultimo" (last-checked) column D
# with new timestamps. Rows 2-43 fall into three 14-row blocks that each
# previously shared one timestamp; the refresh shifts each block to the
# value the block above it had, and stamps the newest block (rows 2-15)
# with the current update time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamps = @{
    2  = 44240.50964031673
    16 = 44240.48844734954
    30 = 44240.46725331018
}

foreach ($startRow in $newTimestamps.Keys) {
    $value = $newTimestamps[$startRow]
    for ($r = $startRow; $r -lt ($startRow + 14); $r++) {
        $ws.Range("D$r").Value = $value
    }
}
